$d = $word.ActiveDocument

# Find the paragraph "Data, Technology and Strategy Consulting" and insert
# three new bullet paragraphs right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Data, Technology and Strategy Consulting") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph 'Data, Technology and Strategy Consulting'"
}

$newLines = @(
    "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters",
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

$insertText = [string]::Join("`r", $newLines) + "`r"

$r = $target.Range
$r.Collapse(0)
$r.InsertAfter($insertText)
